$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("quality_comparison")
$ws2 = $wb.Worksheets.Item("computational_comparison")

# ---------------------------------------------------------------------------
# Build the two border "template" styles once (on sheet1's C1/D1), then
# replicate them onto every other header cell that needs the same frame via
# a format copy/paste, so that the workbook ends up with exactly the two new
# shared styles instead of duplicate/garbage entries for every cell touched.
#   style "topBottom"      -> thin top + thin bottom              (border 4)
#   style "topBottomRight" -> thin top + thin bottom + thin right (border 5)
# ---------------------------------------------------------------------------

$c1 = $ws1.Range("C1")
$c1.ClearFormats()
$c1.Borders.Item(8).LineStyle = 1   # xlEdgeTop
$c1.Borders.Item(8).Weight = 2
$c1.Borders.Item(9).LineStyle = 1   # xlEdgeBottom
$c1.Borders.Item(9).Weight = 2

$d1 = $ws1.Range("D1")
$d1.ClearFormats()
$c1.Copy()
$d1.PasteSpecial(-4122)             # xlPasteFormats
$d1.Borders.Item(10).LineStyle = 1  # xlEdgeRight
$d1.Borders.Item(10).Weight = 2

# Anonymize "fedcore" -> "approach"
$ws1.Range("C2").Value = "approach"

# ---------------------------------------------------------------------------
# Sheet 2: computational_comparison - reuse the templates built above.
# ---------------------------------------------------------------------------
$c1b = $ws2.Range("C1")
$c1b.ClearFormats()
$c1.Copy()
$c1b.PasteSpecial(-4122)            # xlPasteFormats

$d1b = $ws2.Range("D1")
$d1b.ClearFormats()
$d1.Copy()
$d1b.PasteSpecial(-4122)            # xlPasteFormats

$f1b = $ws2.Range("F1")
$f1b.ClearFormats()
$c1.Copy()
$f1b.PasteSpecial(-4122)            # xlPasteFormats

$g1b = $ws2.Range("G1")
$g1b.ClearFormats()
$d1.Copy()
$g1b.PasteSpecial(-4122)            # xlPasteFormats

$excel.CutCopyMode = 0

# Anonymize "fedcore" -> "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# G5 held an empty placeholder cell (0/0 % change) - remove it entirely.
$ws2.Range("G5").ClearContents()

$wb.Save()
